# Fix broken/unvalidated e-mails in the student roster:
#  - row 2 (Aishwarya Singh)   -> Mehul Chopda      (was actually row 3's data)
#  - row 3 (Mehul Chopda)      -> Abhiraj Bishnoi   (was actually row 5's data)
#  - row 4 (Dibyajyoti Ghosh)  -> Aishwarya Singh Bhati (corrected name + new @sitpune.edu.in email)
#  - row 5 (Abhiraj Bishnoi)   -> Akshita Pradhan (new row, new @sitpune.edu.in email)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Mehul Chopda -------------------------------------------------
$ws.Range("A2").Value = 701
$ws.Range("B2").Value = "Mehul Chopda"
$ws.Range("C2").Value = "mehul.chopda@sitpune.edu.in"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 3
$ws.Rows.Item(2).AutoFit()

# --- Row 3: Abhiraj Bishnoi ----------------------------------------------
$ws.Range("A3").Value = 602
$ws.Range("B3").Value = "Abhiraj Bishnoi"
$ws.Range("C3").Value = "abhiraj.bishnoi@sitpune.edu.in"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 1

# --- Row 4: Aishwarya Singh Bhati -----------------------------------------
$ws.Range("A4").Value = 135
$ws.Range("B4").Value = "Aishwarya Singh Bhati"
$ws.Range("C4").Value = "aishwarya.singh@sitpune.edu.in"
$ws.Range("D4").Value = 3
$ws.Range("E4:G4").ClearContents()

# --- Row 5: Akshita Pradhan ------------------------------------------------
$ws.Range("A5").Value = 790
$ws.Range("B5").Value = "Akshita Pradhan"
$ws.Range("C5").Value = "akshita.pradhan@sitpune.edu.in"
$ws.Range("D5").Value = 3
$ws.Range("E5:G5").ClearContents()

# --- Hyperlinks: drop the two stale mailto links, add the two new ones ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:aishwarya.singh@sitpune.edu.in", [type]::Missing, [type]::Missing, "aishwarya.singh@sitpune.edu.in")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:akshita.pradhan@sitpune.edu.in", [type]::Missing, [type]::Missing, "akshita.pradhan@sitpune.edu.in")

# Reuse the already-correctly-styled e-mail cell (C2) format so every
# e-mail cell shares the same hyperlink-coloured font instead of Excel's
# auto-inserted (underlined, themed) "Hyperlink" style. Must happen AFTER
# both the .Value assignments and Hyperlinks.Add above, since either of
# those resets a cell's style.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Selection ends up on A20 in the saved file ---------------------------
$ws.Range("A20").Select()
